$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1918552174508198
$ws.Range("C2").Value = 0.5517837712442845
$ws.Range("D2").Value = 0.5641894466527427
$ws.Range("E2").Value = 0.7511254533383506
$ws.Range("F2").Value = 0.7362272034013972
$ws.Range("G2").Value = 37

$ws.Range("B3").Value = 0.01060547437218578
$ws.Range("C3").Value = 0.4423892562697259
$ws.Range("D3").Value = 0.3511281738231621
$ws.Range("E3").Value = 0.5925606921009544
$ws.Range("F3").Value = 0.600869967357191
$ws.Range("G3").Value = 36

$ws.Range("B4").Value = 0.1023758857763177
$ws.Range("C4").Value = 0.4466085728510511
$ws.Range("D4").Value = 0.3426007271002285
$ws.Range("E4").Value = 0.5853210461791277
$ws.Range("F4").Value = 0.5847120125455174
$ws.Range("G4").Value = 35

$ws.Range("B5").Value = 0.03600675553489324
$ws.Range("C5").Value = 0.4382764901874337
$ws.Range("D5").Value = 0.3537032576543248
$ws.Range("E5").Value = 0.5947295668237159
$ws.Range("F5").Value = 0.602565983338879
$ws.Range("G5").Value = 34

$ws.Range("B6").Value = 0.07600721640812284
$ws.Range("C6").Value = 0.4186182603205892
$ws.Range("D6").Value = 0.3343428317836482
$ws.Range("E6").Value = 0.5782238595765902
$ws.Range("F6").Value = 0.5820939907362127
$ws.Range("G6").Value = 33

$ws.Range("B7").Value = 0.09996734948318312
$ws.Range("C7").Value = 0.3946104424313834
$ws.Range("D7").Value = 0.304156198780297
$ws.Range("E7").Value = 0.5515035800249143
$ws.Range("F7").Value = 0.5510461396922992
$ws.Range("G7").Value = 32

$ws.Range("B8").Value = 0.058982071561956
$ws.Range("C8").Value = 0.4225312785458647
$ws.Range("D8").Value = 0.3413796702787015
$ws.Range("E8").Value = 0.5842770492486433
$ws.Range("F8").Value = 0.59090112965712
$ws.Range("G8").Value = 31

$ws.Range("B9").Value = 0.1107680203475256
$ws.Range("C9").Value = 0.4177884961770141
$ws.Range("D9").Value = 0.3290473224825335
$ws.Range("E9").Value = 0.5736264659885678
$ws.Range("F9").Value = 0.572451866506145
$ws.Range("G9").Value = 30

$ws.Range("B10").Value = 0.07549334707445945
$ws.Range("C10").Value = 0.4371919483332072
$ws.Range("D10").Value = 0.3654170218579554
$ws.Range("E10").Value = 0.6044973299014277
$ws.Range("F10").Value = 0.6103808973489443
$ws.Range("G10").Value = 29

$ws.Range("B11").Value = 0.1033470187209335
$ws.Range("C11").Value = 0.4328287120564761
$ws.Range("D11").Value = 0.3505829942180106
$ws.Range("E11").Value = 0.5921004933438332
$ws.Range("F11").Value = 0.5937098325534103
$ws.Range("G11").Value = 28

